$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "LCL" column (old column B) so everything shifts one column left:
# CBM | LCL | 20' | 40' | 40HC | 45' | Utilization
# -> CBM | 20' | 40' | 40HC | 45' | Utilization
$ws.Columns("B").Delete()

# The "Utilization" column holds percentages formatted as plain text
# (e.g. "52.2%"), not numeric percentages - force text storage so Excel
# doesn't auto-convert the string into a 0.522 percent-formatted number.
$ws.Range("F2:F5").NumberFormat = "@"

# Replace the data rows (2-4) with the new figures, and add a new row 5.
$ws.Range("A2").Value = 29.73
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = "52.2%"

$ws.Range("A3").Value = 11.12
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = "41.2%"

$ws.Range("A4").Value = 41.19
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = "72.3%"

$ws.Range("A5").Value = 32.02
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = "56.2%"

# Drop the number-format override (the original sheet has no explicit
# per-cell styles) now that the text values are locked in.
$ws.Range("F2:F5").ClearFormats()

# Turn on the autofilter for the header row.
$null = $ws.Range("A1").AutoFilter()

# AutoFilter needs the accompanying hidden _FilterDatabase defined name
# (scoped to this sheet) so the filter range is remembered.
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", '=''Sheet''!$A$1')
$filterName.Visible = $false

# Match the saved selection/active cell.
$null = $ws.Range("C4").Select()
